# "Se agrega nueva mercancia e inventario"
# Updates several existing IMEI/serial values in the inventory table and
# appends two new rows (15 and 16) of merchandise/inventory data, plus
# two new notes in column E for rows 12 and 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing cells (rows 9-14) ------------------------------------
$ws.Cells.Item(9, 3).Value = "3052749177"
$ws.Cells.Item(9, 4).Value = "732111193280551"

$ws.Cells.Item(10, 3).Value = "3052754285"
$ws.Cells.Item(10, 4).Value = "732111324709512"

$ws.Cells.Item(11, 3).Value = "3052749177"
$ws.Cells.Item(11, 4).Value = "732111193280551"

$ws.Cells.Item(12, 4).Value = "732111324709673"
# New annotation cell E12 - copy text formatting from an existing "text" cell
$ws.Cells.Item(9, 3).Copy($ws.Cells.Item(12, 5))
$ws.Cells.Item(12, 5).Value = "client nit a nit"

$ws.Cells.Item(13, 4).Value = "732111324709674"
# New annotation cell E13
$ws.Cells.Item(9, 3).Copy($ws.Cells.Item(13, 5))
$ws.Cells.Item(13, 5).Value = "988154393"

# Row 14 - align its formatting with the rest of the block (C/D use the
# same text style as rows 9-13) and update its values
$ws.Cells.Item(9, 3).Copy($ws.Cells.Item(14, 3))
$ws.Cells.Item(14, 3).Value = "3045984556"
$ws.Cells.Item(9, 3).Copy($ws.Cells.Item(14, 4))
$ws.Cells.Item(14, 4).Value = "732111324709675"

# --- Append new rows 15 and 16 --------------------------------------------
$ws.Cells.Item(9, 1).Copy($ws.Cells.Item(15, 1))
$ws.Cells.Item(15, 1).Value = "10960370"
$ws.Cells.Item(9, 3).Copy($ws.Cells.Item(15, 2))
$ws.Cells.Item(15, 2).Value = "36844580"
$ws.Cells.Item(9, 3).Copy($ws.Cells.Item(15, 3))
$ws.Cells.Item(15, 3).Value = "3052754289"
$ws.Cells.Item(9, 3).Copy($ws.Cells.Item(15, 4))
$ws.Cells.Item(15, 4).Value = "732111324709676"

$ws.Cells.Item(9, 1).Copy($ws.Cells.Item(16, 1))
$ws.Cells.Item(16, 1).Value = "10960370"
$ws.Cells.Item(9, 3).Copy($ws.Cells.Item(16, 2))
$ws.Cells.Item(16, 2).Value = "914355426"
$ws.Cells.Item(9, 3).Copy($ws.Cells.Item(16, 3))
$ws.Cells.Item(16, 3).Value = "3046008586"
$ws.Cells.Item(9, 3).Copy($ws.Cells.Item(16, 4))
$ws.Cells.Item(16, 4).Value = "732111193278871"

# --- Update selection to reflect where the user ended up editing ---------
$ws.Range("F16").Select()
